$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Cards"

# Update header row (row 1): shift header meanings, drop the old A1 header
$ws.Range("B1").Value = "card_name"
$ws.Range("C1").Value = "card_cid"
$ws.Range("D1").Value = "card_strats"
$ws.Range("E1").Value = "card_tags"

# Remove the old A1 header cell entirely (cube_name no longer exists)
$ws.Range("A1").Clear()

# Update data row (row 2)
# A2 becomes a numeric 0, formatted like the header cells (copy format from B1)
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0

$ws.Range("B2").Value = "Ornithopter"
$ws.Range("C2").Value = "UR"
$ws.Range("D2").Value = "Arf, Meow"
$ws.Range("E2").Value = "Commander"
